# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): set the new header text first.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header cells (bold, centered, bordered)
# by copying the format from an existing header cell (AC1) onto the new ones.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-52: fill in the team record values (same W/L/T for every row).
$ws.Range("AD2:AD52").Value = 89
$ws.Range("AE2:AE52").Value = 73
$ws.Range("AF2:AF52").Value = 0
